# Updated cryptos list - refreshed prices / 1h volume percentages.
# Leading "'" forces numeric-looking text (e.g. "623.93") to stay as text,
# matching the workbook's existing convention of storing prices as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.060.31"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").Value = "3.584.06"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'623.93"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("D6").Value = "'157.32"
$ws.Range("E6").Value = "  +5.86%  "
$ws.Range("D7").Value = "3.578.75"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  +8.45%  "
$ws.Range("D11").Value = "'7.40"
$ws.Range("E11").Value = "  +7.52%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  +4.54%  "
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  +5.51%  "
$ws.Range("D14").Value = "'33.50"
$ws.Range("E14").Value = "  +6.98%  "
$ws.Range("D15").Value = "4.195.98"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.470.64"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.586.83"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  +5.19%  "
$ws.Range("D20").Value = "'16.09"
$ws.Range("E20").Value = "  +7.09%  "
$ws.Range("E21").Value = "  +12.52%  "
$ws.Range("D22").Value = "'461.27"
$ws.Range("E22").Value = "  +3.86%  "
$ws.Range("D23").Value = "'0.643"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'78.78"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").Value = "'0.0000134"
$ws.Range("E25").Value = "  +8.12%  "
$ws.Range("D26").Value = "'10.66"
$ws.Range("E26").Value = "  +5.16%  "
$ws.Range("D27").Value = "3.733.73"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'9.17"
$ws.Range("E29").Value = "  +10.44%  "
$ws.Range("D30").Value = "'2.63"
$ws.Range("E30").Value = "  +3.91%  "
$ws.Range("E31").Value = "  +8.46%  "
$ws.Range("D32").Value = "'0.172"
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'6.51"
$ws.Range("E34").Value = "  +6.31%  "
$ws.Range("D35").Value = "'26.41"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'1.93"
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("D37").Value = "3.583.55"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("D38").Value = "'8.39"
$ws.Range("E38").Value = "  +5.17%  "
$ws.Range("E39").Value = "  +8.85%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'179.82"
$ws.Range("E41").Value = "  +5.65%  "
$ws.Range("D42").Value = "'0.0922"
$ws.Range("E42").Value = "  +6.15%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'5.69"
$ws.Range("E44").Value = "  +4.64%  "
$ws.Range("D45").Value = "'31.32"
$ws.Range("E45").Value = "  +20.46%  "
$ws.Range("D46").Value = "'0.906"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "'1.37"
$ws.Range("E47").Value = "  +10.62%  "
$ws.Range("D48").Value = "'45.93"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").Value = "'2.73"
$ws.Range("E49").Value = "  +8.71%  "
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").Value = "'0.266"
$ws.Range("E51").Value = "  +8.30%  "
